$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.063.96"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.018.57"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0790"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "2.315.32"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.744"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "2.035.21"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "37.057.84"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.126"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").Value = "1.477.26"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("E45").Value = "  -5.08%  "
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "2.204.44"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.07%  "
